# Insert a new glossary entry "finely-ground" right after the existing
# "finely dispersed" entry (and before "flat-lying"), matching the stylesheet
# formatting used by neighbouring entries (Normal style, non-italic run).

$d = $word.ActiveDocument

# Locate the paragraph that contains the "finely dispersed" entry.
$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($text -eq "finely dispersed") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'finely dispersed' paragraph"
}

# Work on the paragraph's range, collapsed to its end (after the text but
# before the paragraph mark) so a new paragraph can be appended right after it.
$insertionRange = $target.Range
$insertionRange.Collapse(0)  # wdCollapseEnd

$insertionRange.InsertParagraphAfter()

# Move into the freshly created paragraph (the one following our target).
$newPara = $target.Next()
$newRange = $newPara.Range
$newRange.Collapse(1)  # wdCollapseStart -- start of the new, empty paragraph

# Insert the text first, then format only the inserted run (not the
# paragraph mark), so the new paragraph's own pPr/rPr stays empty -- just
# like the other plain (non-tracked) entries in this glossary.
$newRange.InsertAfter("finely-ground")
$newRange.Font.Italic = $false
